# Auto-generated update script for cryptos.xlsx
# Applies per-cell text/value updates matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cells whose new text could be mis-typed as a number by Excel's
#     automatic type inference (e.g. "601.11", "1.00") get a quick
#     "@" (Text) format applied before the write, then the format is
#     reset back to the default "Normal" style so no stray formatting
#     is left behind, matching the source (which used the default style).

$textCells = @('D5', 'D6', 'D8', 'D10', 'D12', 'D14', 'D19', 'D20', 'D23', 'D25', 'D27', 'D30', 'D32', 'D33', 'D35', 'D36', 'D37', 'D38', 'D39', 'D40', 'D41', 'D43', 'D46', 'D47', 'D48', 'D49', 'D50')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D5').Value = '601.11'
$ws.Range('D6').Value = '155.56'
$ws.Range('D8').Value = '0.547'
$ws.Range('D10').Value = '0.139'
$ws.Range('D12').Value = '5.24'
$ws.Range('D14').Value = '27.98'
$ws.Range('D19').Value = '11.40'
$ws.Range('D20').Value = '365.32'
$ws.Range('D23').Value = '4.93'
$ws.Range('D25').Value = '72.97'
$ws.Range('D27').Value = '10.11'
$ws.Range('D30').Value = '584.45'
$ws.Range('D32').Value = '1.42'
$ws.Range('D33').Value = '8.02'
$ws.Range('D35').Value = '0.132'
$ws.Range('D36').Value = '1.00'
$ws.Range('D37').Value = '1.55'
$ws.Range('D38').Value = '160.47'
$ws.Range('D39').Value = '1.93'
$ws.Range('D40').Value = '19.33'
$ws.Range('D41').Value = '5.40'
$ws.Range('D43').Value = '2.66'
$ws.Range('D46').Value = '40.74'
$ws.Range('D47').Value = '1.00'
$ws.Range('D48').Value = '156.28'
$ws.Range('D49').Value = '3.74'
$ws.Range('D50').Value = '22.07'

foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}

# --- Remaining cells (coin names, links, percentages) are safe to set
#     directly; Excel keeps them as plain text.

$ws.Range('D2').Value = '68.792.21'
$ws.Range('E2').Value = '  +2.14%  '
$ws.Range('D3').Value = '2.649.55'
$ws.Range('E3').Value = '  +1.45%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('E5').Value = '  +1.67%  '
$ws.Range('E6').Value = '  +3.31%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('E8').Value = '  -0.78%  '
$ws.Range('D9').Value = '2.647.74'
$ws.Range('E9').Value = '  +1.60%  '
$ws.Range('E10').Value = '  +13.37%  '
$ws.Range('E11').Value = '  -0.34%  '
$ws.Range('E12').Value = '  +0.97%  '
$ws.Range('E13').Value = '  +1.85%  '
$ws.Range('E14').Value = '  +2.02%  '
$ws.Range('E15').Value = '  +6.21%  '
$ws.Range('D16').Value = '3.132.38'
$ws.Range('E16').Value = '  +1.34%  '
$ws.Range('D17').Value = '68.732.42'
$ws.Range('E17').Value = '  +2.14%  '
$ws.Range('D18').Value = '2.638.12'
$ws.Range('E18').Value = '  +0.86%  '
$ws.Range('E19').Value = '  +3.30%  '
$ws.Range('E20').Value = '  -0.55%  '
$ws.Range('E21').Value = '  +1.67%  '
$ws.Range('E22').Value = '  +0.00%  '
$ws.Range('E23').Value = '  +2.16%  '
$ws.Range('E24').Value = '  +4.54%  '
$ws.Range('E25').Value = '  +10.01%  '
$ws.Range('E26').Value = '  +0.06%  '
$ws.Range('E27').Value = '  +1.06%  '
$ws.Range('E28').Value = '  +6.78%  '
$ws.Range('D29').Value = '2.779.66'
$ws.Range('E29').Value = '  +1.31%  '
$ws.Range('E30').Value = '  +0.71%  '
$ws.Range('E31').Value = '  +0.00%  '
$ws.Range('E32').Value = '  +3.76%  '
$ws.Range('E33').Value = '  +4.88%  '
$ws.Range('E34').Value = '  +3.20%  '
$ws.Range('E35').Value = '  +4.95%  '
$ws.Range('E36').Value = '  -0.05%  '
$ws.Range('E38').Value = '  +3.30%  '
$ws.Range('B39').Value = 'Stacks'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('E39').Value = '  +3.75%  '
$ws.Range('B40').Value = 'EthereumClassic'
$ws.Range('C40').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('E40').Value = '  +1.72%  '
$ws.Range('E41').Value = '  +3.44%  '
$ws.Range('E42').Value = '  +0.73%  '
$ws.Range('E44').Value = '  +5.36%  '
$ws.Range('D45').Value = '0.0₆0322'
$ws.Range('E45').Value = '  +10.76%  '
$ws.Range('B46').Value = 'OKB'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('E46').Value = '  -0.10%  '
$ws.Range('B47').Value = 'USDe'
$ws.Range('C47').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('E47').Value = '  +0.03%  '
$ws.Range('E48').Value = '  +0.44%  '
$ws.Range('E49').Value = '  +0.54%  '
$ws.Range('E50').Value = '  +2.95%  '
$ws.Range('E51').Value = '  +1.31%  '
